# corrigindo ordenação dos campos
#
# Update the rebalancing inputs on Plan1: the cash available to invest (G3),
# the current unit prices for each position (I6:I11), and flag two positions
# (EMBR3 / GGBR4, rows 7 and 8) as "fixed" (T7/T8 = "X") so their target
# allocation is pinned instead of recomputed. All of the other cells in the
# sheet are formulas and simply recalculate from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cash on hand (feeds G2 = R12 + G3)
$ws.Range("G3").Value = 2500

# Current quoted price per share for each row (I6:I11)
$ws.Range("I6").Value  = 44.74
$ws.Range("I7").Value  = 22.79
$ws.Range("I8").Value  = 15.46
$ws.Range("I9").Value  = 10.35
$ws.Range("I10").Value = 52.76
$ws.Range("I11").Value = 17.05

# Mark EMBR3 (row 7) and GGBR4 (row 8) as fixed/locked positions
$ws.Range("T7").Value = "X"
$ws.Range("T8").Value = "X"

# Reflect the user's last selection/scroll position on the sheet
[void]$ws.Range("N11").Select()
